# Scheduled-runner refresh of market/profit figures across the Leve
# profitability sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N)
# for the rows whose underlying market data changed, and clears the
# LeveProfitHQ (N) / LeveProfitNQ (M) cells that no longer apply where
# the HQ price data dropped out (no HQ list price -> no HQ profit figure).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 100
$ws.Range("I5").Value = 100
$ws.Range("K5").Value = 100
$ws.Range("M5").Value = 15

$ws.Range("H53").Value = 346.91666
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()

$ws.Range("H88").Value = 991.5
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value = 991.5
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws.Range("H131").Value = 9047.5
$ws.Range("J131").Value = 8000
$ws.Range("L131").Value = 24000
$ws.Range("N131").Value = -34080

$ws.Range("H138").Value = 2729.8794
$ws.Range("J138").Value = 3041.932
$ws.Range("L138").Value = 9125.795999999998
$ws.Range("N138").Value = -19405.796

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1040.0769
$ws.Range("I2").Value = 1184.3636
$ws.Range("J2").Value = 246.5
$ws.Range("K2").Value = 1184.3636
$ws.Range("L2").Value = 246.5
$ws.Range("M2").Value = -1071.3636
$ws.Range("N2").Value = -472.5

$ws.Range("H32").Value = 8211.714
$ws.Range("I32").Value = 8211.714
$ws.Range("K32").Value = 8211.714
$ws.Range("M32").Value = -7924.714

$ws.Range("H61").Value = 2223.7144
$ws.Range("I61").Value = 1579.9375
$ws.Range("K61").Value = 1579.9375
$ws.Range("M61").Value = -1367.9375

$ws.Range("H116").Value = 1040.0769
$ws.Range("I116").Value = 1184.3636
$ws.Range("J116").Value = 246.5
$ws.Range("K116").Value = 1184.3636
$ws.Range("L116").Value = 246.5
$ws.Range("M116").Value = 1109.6364
$ws.Range("N116").Value = -4834.5

$ws.Range("H136").Value = 2223.7144
$ws.Range("I136").Value = 1579.9375
$ws.Range("K136").Value = 4739.8125
$ws.Range("M136").Value = -2189.8125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1040.0769
$ws.Range("I3").Value = 1184.3636
$ws.Range("J3").Value = 246.5
$ws.Range("K3").Value = 1184.3636
$ws.Range("L3").Value = 246.5
$ws.Range("M3").Value = -1070.3636
$ws.Range("N3").Value = -474.5

$ws.Range("H86").Value = 5099.8
$ws.Range("I86").Value = 5124.75
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 5124.75
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -4001.75
$ws.Range("N86").Value = -7246

$ws.Range("H89").Value = 5099.8
$ws.Range("I89").Value = 5124.75
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 25623.75
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -20007.75
$ws.Range("N89").Value = -36232

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 436.75
$ws.Range("J7").Value = 733.3333
$ws.Range("L7").Value = 733.3333
$ws.Range("N7").Value = -959.3333

$ws.Range("H16").Value = 1281.625
$ws.Range("I16").Value = 1378
$ws.Range("J16").Value = 992.5
$ws.Range("K16").Value = 1378
$ws.Range("L16").Value = 992.5
$ws.Range("M16").Value = -1091
$ws.Range("N16").Value = -1566.5

$ws.Range("H58").Value = 2719.4443
$ws.Range("I58").Value = 2262.3333
$ws.Range("J58").Value = 2948
$ws.Range("K58").Value = 2262.3333
$ws.Range("L58").Value = 2948
$ws.Range("M58").Value = -2059.3333
$ws.Range("N58").Value = -3354

$ws.Range("H105").Value = 1710.75
$ws.Range("I105").Value = 1769.4286
$ws.Range("K105").Value = 1769.4286
$ws.Range("M105").Value = -22.42859999999996

$ws.Range("H113").Value = 1281.625
$ws.Range("I113").Value = 1378
$ws.Range("J113").Value = 992.5
$ws.Range("K113").Value = 1378
$ws.Range("L113").Value = 992.5
$ws.Range("M113").Value = 792
$ws.Range("N113").Value = -5332.5

$ws.Range("H136").Value = 2719.4443
$ws.Range("I136").Value = 2262.3333
$ws.Range("J136").Value = 2948
$ws.Range("K136").Value = 6786.999899999999
$ws.Range("L136").Value = 8844
$ws.Range("M136").Value = -4236.999899999999
$ws.Range("N136").Value = -13944

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 2966.6667
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()

$ws.Range("H55").Value = 30303992
$ws.Range("I55").Value = 30303992
$ws.Range("K55").Value = 90911976
$ws.Range("M55").Value = -90911799

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3053.7334
$ws.Range("I132").Value = 1274.3334
$ws.Range("J132").Value = 3498.5833
$ws.Range("K132").Value = 3823.0002
$ws.Range("L132").Value = 10495.7499
$ws.Range("M132").Value = -1293.0002
$ws.Range("N132").Value = -15555.7499

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1624.5
$ws.Range("I46").Value = 1624.5
$ws.Range("K46").Value = 1624.5
$ws.Range("M46").Value = -1436.5

$ws.Range("H61").Value = 3665.6667
$ws.Range("I61").Value = 2999
$ws.Range("K61").Value = 2999
$ws.Range("M61").Value = -2797

$ws.Range("H68").Value = 2894
$ws.Range("I68").Value = 2882.125
$ws.Range("J68").Value = 2989
$ws.Range("K68").Value = 2882.125
$ws.Range("L68").Value = 2989
$ws.Range("M68").Value = -2133.125
$ws.Range("N68").Value = -4487

$ws.Range("H71").Value = 2894
$ws.Range("I71").Value = 2882.125
$ws.Range("J71").Value = 2989
$ws.Range("K71").Value = 14410.625
$ws.Range("L71").Value = 14945
$ws.Range("M71").Value = -10666.625
$ws.Range("N71").Value = -22433

$ws.Range("H113").Value = 3665.6667
$ws.Range("I113").Value = 2999
$ws.Range("K113").Value = 2999
$ws.Range("M113").Value = -829

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").ClearContents()

$ws.Range("H107").Value = 486.33334
$ws.Range("I107").Value = 486.33334
$ws.Range("K107").Value = 1459.00002
$ws.Range("M107").Value = 460.9999800000001
